$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 (D30 = "?") becomes row 31, with two new cells added alongside it:
# A31 = 20-03-2017 (date serial 42814), same date-format style as A29/A28
# B31 = "Tolga Yasin Kücük" (same value as B29/B28)
# D31 keeps the same value that used to be in D30 ("?")

$oldD30 = $ws.Range("D30").Value2
$oldB29 = $ws.Range("B29").Value2

$ws.Range("D30").ClearContents()

$ws.Range("A31").Value = 42814
$ws.Range("A31").NumberFormat = "m/d/yy"
$ws.Range("B31").Value = $oldB29
$ws.Range("D31").Value = $oldD30

[void]$ws.Range("D31").Select()

# Best-effort: mirror the author's scroll position (topLeftCell A25) in the
# active window's scroll state.
try {
    $win = $excel.ActiveWindow
    $win.ScrollColumn = 1
    $win.ScrollRow = 25
} catch {
}
